# ExpenseSheet2017-2018.xlsx — remove the second expense entry (row 3):
#   Sl.No 2 / 20-Jul-2017 / "Knowtefy DSC & DIN" / "Rabindra (CA)" /
#   4000 / 0 / 4000 / "Cash Payment"
# is deleted from the sheet, leaving row 3 blank (matching the already
# blank rows 4-9), and the active selection moves to H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 currently carries the date-formatted style (s=2) on B3 because it
# held a date value; once the row is cleared it should look like the other
# empty rows below it (all cells style s=1). Copy B4's (already-blank,
# border-only) formatting onto B3 before wiping the values so the cell
# drops its number format along with its content.
$ws.Range("B4").Copy()
$ws.Range("B3").PasteSpecial(-4122)  # xlPasteFormats

# Clear out the entry itself.
$ws.Range("A3:H3").ClearContents()

# Leave the selection where the user last clicked.
$ws.Range("H2").Select()
